$wb = $excel.ActiveWorkbook

$langs = @(
    @{ Sheet = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackTime = "2016-01-27 07:49:28" },
    @{ Sheet = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackTime = "2016-01-27 07:49:46" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    foreach ($row in 2, 3) {
        # Status: handback is complete and in sync with en-US source
        $ws.Cells.Item($row, 2).Value = "Handed back: in sync with en-US"

        # E = Latest Target File, F = Latest Handback File: populated now that the
        # file has been handed back
        $ws.Cells.Item($row, 5).Value = "a.md"
        $ws.Cells.Item($row, 5).Style = $ws.Cells.Item($row, 1).Style

        $ws.Cells.Item($row, 6).Value = $lang.Xlf
        $ws.Cells.Item($row, 6).Style = $ws.Cells.Item($row, 3).Style

        # G = Latest Handback DateTime: stamp the real handback time
        $ws.Cells.Item($row, 7).Value = $lang.HandbackTime
    }
}
